# WIP: Improve images, and practice rounds UX.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Unfreeze panes (sheetView previously had a frozen B2 pane) ---
$excel.ActiveWindow.FreezePanes = $false

# --- 2. Widen column D to fit the new, longer "correct" answer text ---
$ws.Columns(4).ColumnWidth = 29.1484

# --- 3. Update wording of the "correct" values ---
# "both" -> "left,right" for the existing rows (D2:D6)
$ws.Range("D2:D6").Value = "left,right"

# Fill in the previously-blank "correct" cells for the new practice rows
# (rows 31-36, where left_color == right_color) with "left,right,"
$ws.Range("D31:D36").Value = "left,right,"

# --- 4. Column D formatting: left-align + slightly smaller font (11pt) ---
$dRange = $ws.Range("D1:D36")
$dRange.HorizontalAlignment = -4131   # xlHAlignLeft
$dRange.Font.Size = 11
$dRange.Font.Bold = $false

# --- 5. Give the data cells (everything except the left-most "group" column
#         and the header row) an explicit new background fill ---
$ws.Range("B2:G36").Interior.ColorIndex = 6
